$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Doing")
$ws2 = $wb.Worksheets.Item("Done")

# 1. Update T6 on "Doing" sheet with a new lesson date for 刘桐语
$ws1.Range("T6").Value = 45703

# 2. Insert a new row at row 13 on "Doing" sheet (shifts rows 19-30 down to 20-31)
$ws1.Rows("13:13").Insert()

# 3. Populate the new row 13 with the 张佳宁 tutoring record (moved back from "Done")
$ws1.Range("A13").Value = "初三"
$ws1.Range("B13").Value = "物理"
$ws1.Range("C13").Value = "一尔优教育"
$ws1.Range("D13").Value = "240/2h"
$ws1.Range("E13").Value = "周六12：00-14：00"
$ws1.Range("F13").Value = "张佳宁"
$ws1.Range("G13").Value = "耳鼻喉科专科医院旁巷子201室"
$ws1.Range("H13").Value = "13中"
$ws1.Range("I13").Value = "85/100"
$ws1.Range("J13").Value = "/"
$ws1.Range("K13").Value = 45451
$ws1.Range("L13").Value = 45460
$ws1.Range("M13").Value = 45465
$ws1.Range("N13").Value = 45514
$ws1.Range("O13").Value = 45527
$ws1.Range("P13").Value = 45658
$ws1.Range("Q13").Value = 45703

# Match styles used for the row (A-F: s5, G-I: s4, J-O: s12, P: s9, Q: s26, R-T: s16)
$ws1.Range("A13:F13").Style = $ws1.Range("A11").Style
$ws1.Range("G13:I13").Style = $ws1.Range("G11").Style
$ws1.Range("J13:O13").Style = $ws1.Range("J9").Style
$ws1.Range("P13").Style = $ws1.Range("O2").Style
$ws1.Range("Q13").Style = $ws1.Range("S3").Style
$ws1.Range("R13:T13").Style = $ws1.Range("P8").Style

# 4. Clear the old record row (row 27) on "Done" sheet, since it moved to "Doing"
$ws2.Range("A27:T27").ClearContents()
